$wb = $excel.ActiveWorkbook

# Insert the new "practiceQ" sheet between "textEditor" and "credentials".
$textSheet = $wb.Worksheets.Item("textEditor")
$newSheet = $wb.Worksheets.Add($null, $textSheet)
$newSheet.Name = "practiceQ"

# Header row.
$newSheet.Range("A1").Value = "Valid code"
$newSheet.Range("B1").Value = "Invalid Code"

# Data row (set B2 before A2 so shared-string indices come out abc=15, print'hello'=16).
$newSheet.Range("B2").Value = "abc"
$newSheet.Range("A2").Value = "print 'hello'"

# Formatting.
$newSheet.Columns.Item(1).ColumnWidth = 68.14
$newSheet.Range("A2").WrapText = $true
$newSheet.Range("B2").VerticalAlignment = -4160
$newSheet.Range("B7").VerticalAlignment = -4160
$newSheet.Rows.Item(2).RowHeight = 16

# Selection/active cell as left by the author.
[void]$newSheet.Range("A6").Select()
